$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

$ws.Range("G4").Value = "BEISHUI_DESC"

$ws.Range("H5").Select()
